$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.777.46'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.91%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.909.98'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.39%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.09%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.36'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.47%  '

$ws.Range('E6').Value = '  -0.05%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4748'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.65%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2871'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.05%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06679'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.00%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.80'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.87%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '101.29'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.51%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07798'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.78%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.928.64'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.62%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.205'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.36%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6857'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.93%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '287.66'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +9.58%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '30.867.76'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.14%  '

$ws.Range('B18').Value = 'Dai'
$ws.Range('C18').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.002'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.14%  '

$ws.Range('B19').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C19').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '2.178.24'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.47%  '

$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007519'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.30%  '

$ws.Range('B21').Value = 'Avalanche'
$ws.Range('C21').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.79'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.58%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.444'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.53%  '

$ws.Range('B23').Value = 'BitDAO'
$ws.Range('C23').Value = 'https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.4664'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.62%  '

$ws.Range('B24').Value = 'BinanceUSD'
$ws.Range('C24').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.002'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.20%  '

$ws.Range('B25').Value = 'Chainlink'
$ws.Range('C25').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '6.309'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.46%  '

$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.399'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.15%  '

$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '168.03'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.03%  '

$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.48'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.60%  '

$ws.Range('B29').Value = 'LidoDAOToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.049'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.33%  '

$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.385'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.13%  '

$ws.Range('B31').Value = 'Stellar'
$ws.Range('C31').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09945'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.96%  '

$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.504'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.78%  '

$ws.Range('B33').Value = 'PancakeSwap'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.522'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.39%  '

$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.267'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.80%  '

$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.04780'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.44%  '

$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7233'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.27%  '

$ws.Range('B37').Value = 'ARBITRUM'
$ws.Range('C37').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.117'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.34%  '

$ws.Range('B38').Value = 'HuobiToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.724'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.58%  '

$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01918'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.51%  '

$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.839'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +8.47%  '

$ws.Range('B41').Value = 'MXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.619'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.77%  '

$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '74.35'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.80%  '

$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.997'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.15%  '

$ws.Range('B44').Value = 'TrustWalletToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.8640'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.28%  '

$ws.Range('B45').Value = 'TheSandbox'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4293'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.47%  '

$ws.Range('B46').Value = 'Quant'
$ws.Range('C46').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '105.01'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.80%  '

$ws.Range('B47').Value = 'PaxDollar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.002'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.08%  '

$ws.Range('B48').Value = 'Maker'
$ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.003.68'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.55%  '

$ws.Range('B49').Value = 'Aptos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.425'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.03%  '

$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.239'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.27%  '

$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.1188'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.15%  '
